$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.071.62'
$ws.Range('E2').Value = '  -0.44%  '
$ws.Range('D3').Value = '1.874.00'
$ws.Range('E3').Value = '  -1.87%  '
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '313.58'
$ws.Range('E5').Value = '  -0.33%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  +0.18%  '
$ws.Range('E7').Value = '  -0.82%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3839'
$ws.Range('E8').Value = '  -2.31%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.08625'
$ws.Range('E9').Value = '  -7.63%  '
$ws.Range('E10').Value = '  -2.20%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '41.57'
$ws.Range('E11').Value = '  -0.83%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '6.306'
$ws.Range('E12').Value = '  -1.54%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '20.64'
$ws.Range('E13').Value = '  -1.20%  '
$ws.Range('D14').Value = '1.878.00'
$ws.Range('E14').Value = '  -1.62%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.190'
$ws.Range('E15').Value = '  -1.83%  '
$ws.Range('E16').Value = '  +0.19%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.00001097'
$ws.Range('E17').Value = '  -2.19%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '90.84'
$ws.Range('E18').Value = '  -1.70%  '
$ws.Range('E19').Value = '  +0.16%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '18.04'
$ws.Range('E20').Value = '  +0.35%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.077'
$ws.Range('E22').Value = '  -2.50%  '
$ws.Range('D23').Value = '28.108.94'
$ws.Range('E23').Value = '  -0.51%  '
$ws.Range('E24').Value = '  -0.86%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.266'
$ws.Range('E25').Value = '  -2.50%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.578'
$ws.Range('E26').Value = '  -0.69%  '
$ws.Range('D27').Value = '2.090.66'
$ws.Range('E27').Value = '  -1.43%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '20.72'
$ws.Range('E28').Value = '  -1.89%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '126.25'
$ws.Range('E30').Value = '  -0.75%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.1053'
$ws.Range('E31').Value = '  -2.02%  '
$ws.Range('E32').Value = '  -4.28%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '5.602'
$ws.Range('E33').Value = '  -0.85%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.593'
$ws.Range('E34').Value = '  -0.57%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '9.608'
$ws.Range('E35').Value = '  -1.18%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.02438'
$ws.Range('E36').Value = '  +0.51%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.06572'
$ws.Range('E37').Value = '  -1.74%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.2173'
$ws.Range('E38').Value = '  -1.61%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.209'
$ws.Range('E39').Value = '  -2.91%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.237'
$ws.Range('E40').Value = '  -3.66%  '
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.6374'
$ws.Range('E41').Value = '  -2.58%  '
$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '11.50'
$ws.Range('E42').Value = '  -0.31%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '4.888'
$ws.Range('E43').Value = '  -2.58%  '
$ws.Range('B44').Value = 'Decentraland'
$ws.Range('C44').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.5996'
$ws.Range('E44').Value = '  -2.05%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '13.13'
$ws.Range('E45').Value = '  -1.44%  '
$ws.Range('E46').Value = '  -0.29%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.674'
$ws.Range('E47').Value = '  -1.43%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.226'
$ws.Range('E48').Value = '  +3.14%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.989'
$ws.Range('E49').Value = '  -1.67%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '121.43'
$ws.Range('E50').Value = '  -1.43%  '
$ws.Range('E51').Value = '  +2.13%  '
